$d = $word.ActiveDocument

# New (Slovak) sentence that replaces every paragraph currently reading
# "V roku 2018 môžete pozorovať súhvezdie Perseus: 30. októbra - 8. novembra
#  a 29. novembra - 8. decembra" (regardless of how many runs / trailing
# whitespace runs it is currently split across). The replacement collapses
# the whole paragraph body down to a single run carrying no direct run
# formatting (no <w:rPr>), per the target markup.
$newText = "V roku Pegasus: 8.-17. októbra, 7.-16. novembra,"

$pkgTemplate = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>{0}</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Forward-only helper: this runtime's LastIndexOf(value, startIndex) does not
# honour the start-index bound, so locate the paragraph's opening tag by
# scanning forward from $fromIdx instead of scanning backward from the end.
function FindParaStart($xml, $fromIdx) {
    $a = $xml.IndexOf("<w:p ", $fromIdx)
    $b = $xml.IndexOf("<w:p>", $fromIdx)
    if ($a -eq -1) { return $b }
    if ($b -eq -1) { return $a }
    if ($a -lt $b) { return $a } else { return $b }
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $r = $para.Range
    $t = $r.Text

    if ($t.Contains("V roku 2018")) {
        # Pull this paragraph's own OOXML so we can keep its <w:p ...> opening
        # tag / attributes and its <w:pPr> (paragraph formatting) untouched -
        # only the run content inside the paragraph needs to change.
        $xml = $r.WordOpenXML
        $bodyIdx = $xml.IndexOf("<w:body>")
        $pIdx = FindParaStart $xml $bodyIdx
        $pPrEndIdx = $xml.IndexOf("</w:pPr>", $pIdx)
        $pPrFull = $xml.Substring($pIdx, $pPrEndIdx + 8 - $pIdx)

        $newRun = "<w:r><w:t>" + $newText + "</w:t></w:r>"
        $newParaXml = $pPrFull + $newRun + "</w:p>"

        $pkg = [string]::Format($pkgTemplate, $newParaXml)

        $r.InsertXML($pkg)
    }
}
